{"js": "// Replace three-digit-division-answer text runs throughout the document.\n// Each target string is unique within the document, so a simple\n// search-and-replace (matching the whole run text) is sufficient and\n// will not disturb any other content or formatting.\nconst replacements = [\n  [\"283\\u00F77=40, 3\", \"664\\u00F72=332, 0\"],\n  [\"550\\u00F74=137, 2\", \"168\\u00F79=18, 6\"],\n  [\"377\\u00F74=94, 1\", \"352\\u00F73=117, 1\"],\n  [\"167\\u00F72=83, 1\", \"633\\u00F75=126, 3\"],\n  [\"312\\u00F77=44, 4\", \"239\\u00F78=29, 7\"],\n  [\"975\\u00F73=325, 0\", \"517\\u00F75=103, 2\"],\n  [\"894\\u00F79=99, 3\", \"211\\u00F76=35, 1\"],\n  [\"876\\u00F79=97, 3\", \"668\\u00F76=111, 2\"],\n  [\"581\\u00F77=83, 0\", \"785\\u00F73=261, 2\"],\n  [\"661\\u00F78=82, 5\", \"726\\u00F72=363, 0\"],\n  [\"307\\u00F72=153, 1\", \"548\\u00F78=68, 4\"],\n  [\"988\\u00F73=329, 1\", \"573\\u00F74=143, 1\"],\n  [\"575\\u00F76=95, 5\", \"880\\u00F76=146, 4\"],\n  [\"473\\u00F74=118, 1\", \"170\\u00F74=42, 2\"],\n  [\"173\\u00F79=19, 2\", \"537\\u00F73=179, 0\"],\n  [\"792\\u00F78=99, 0\", \"433\\u00F72=216, 1\"],\n  [\"140\\u00F77=20, 0\", \"206\\u00F77=29, 3\"],\n  [\"969\\u00F74=242, 1\", \"850\\u00F78=106, 2\"],\n  [\"784\\u00F79=87, 1\", \"180\\u00F76=30, 0\"],\n  [\"921\\u00F78=115, 1\", \"409\\u00F75=81, 4\"],\n  [\"230\\u00F72=115, 0\", \"115\\u00F73=38, 1\"],\n  [\"284\\u00F79=31, 5\", \"404\\u00F72=202, 0\"],\n  [\"308\\u00F76=51, 2\", \"533\\u00F77=76, 1\"],\n  [\"921\\u00F72=460, 1\", \"125\\u00F76=20, 5\"],\n  [\"178\\u00F77=25, 3\", \"124\\u00F77=17, 5\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace three-digit-division-answer text throughout the document.\n# Each target string occurs exactly once in the document, so a plain\n# Find/Replace (whole match, case-sensitive) for each pair is sufficient\n# and leaves every other run/formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"283\u00f77=40, 3\", \"664\u00f72=332, 0\"),\n    @(\"550\u00f74=137, 2\", \"168\u00f79=18, 6\"),\n    @(\"377\u00f74=94, 1\", \"352\u00f73=117, 1\"),\n    @(\"167\u00f72=83, 1\", \"633\u00f75=126, 3\"),\n    @(\"312\u00f77=44, 4\", \"239\u00f78=29, 7\"),\n    @(\"975\u00f73=325, 0\", \"517\u00f75=103, 2\"),\n    @(\"894\u00f79=99, 3\", \"211\u00f76=35, 1\"),\n    @(\"876\u00f79=97, 3\", \"668\u00f76=111, 2\"),\n    @(\"581\u00f77=83, 0\", \"785\u00f73=261, 2\"),\n    @(\"661\u00f78=82, 5\", \"726\u00f72=363, 0\"),\n    @(\"307\u00f72=153, 1\", \"548\u00f78=68, 4\"),\n    @(\"988\u00f73=329, 1\", \"573\u00f74=143, 1\"),\n    @(\"575\u00f76=95, 5\", \"880\u00f76=146, 4\"),\n    @(\"473\u00f74=118, 1\", \"170\u00f74=42, 2\"),\n    @(\"173\u00f79=19, 2\", \"537\u00f73=179, 0\"),\n    @(\"792\u00f78=99, 0\", \"433\u00f72=216, 1\"),\n    @(\"140\u00f77=20, 0\", \"206\u00f77=29, 3\"),\n    @(\"969\u00f74=242, 1\", \"850\u00f78=106, 2\"),\n    @(\"784\u00f79=87, 1\", \"180\u00f76=30, 0\"),\n    @(\"921\u00f78=115, 1\", \"409\u00f75=81, 4\"),\n    @(\"230\u00f72=115, 0\", \"115\u00f73=38, 1\"),\n    @(\"284\u00f79=31, 5\", \"404\u00f72=202, 0\"),\n    @(\"308\u00f76=51, 2\", \"533\u00f77=76, 1\"),\n    @(\"921\u00f72=460, 1\", \"125\u00f76=20, 5\"),\n    @(\"178\u00f77=25, 3\", \"124\u00f77=17, 5\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
